# Apply updated metric values to the keras_metrics workbook.
# Sheet 1: "Validation Metrics"
# Sheet 2: "Classification Report"
# Sheet 3: "Confusion Matrix"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Validation Metrics ---
$ws1 = $wb.Worksheets.Item("Validation Metrics")
$ws1.Range("B2").Value = 0.2823184728622437
$ws1.Range("B3").Value = 0.9249129891395569

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 3 (class "1")
$ws2.Range("B3").Value = 0.9303879310344828
$ws2.Range("C3").Value = 0.9853914631362702
$ws2.Range("D3").Value = 0.9571000997672099
$ws2.Range("E3").Value = 4381

# Row 5 (accuracy)
$ws2.Range("B5").Value = 0.9249130089738111
$ws2.Range("C5").Value = 0.9249130089738111
$ws2.Range("D5").Value = 0.9249130089738111
$ws2.Range("E5").Value = 0.9249130089738111

# Row 6 (macro avg)
$ws2.Range("B6").Value = 0.8756221996196324
$ws2.Range("C6").Value = 0.95093048771209
$ws2.Range("D6").Value = 0.9047301426895916
$ws2.Range("E6").Value = 16381

# Row 7 (weighted avg)
$ws2.Range("B7").Value = 0.9391520998437554
$ws2.Range("C7").Value = 0.9249130089738111
$ws2.Range("D7").Value = 0.9274336582864814
$ws2.Range("E7").Value = 16381

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("C3").Value = 4317
